$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 135 in place: date 44400 -> 44448, volumen (J) 2000 -> 1000 ---
$ws.Cells.Item(135, 4).Value = 44448   # D135 Fecha
$ws.Cells.Item(135, 10).Value = 1000   # J135 Volumen

# --- Update row 136 in place: date 44400 -> 44448, volumen (J) 1000 -> 500 ---
$ws.Cells.Item(136, 4).Value = 44448   # D136 Fecha
$ws.Cells.Item(136, 10).Value = 500    # J136 Volumen

# --- Append new row 137 (copy of prior row135 data: "Primera", Volumen 2000) ---
$ws.Cells.Item(137, 1).Value = 11
$ws.Cells.Item(137, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(137, 3).Value = "Bíobío"
$ws.Cells.Item(137, 4).Value = 44400
$ws.Cells.Item(137, 5).Value = 8
$ws.Cells.Item(137, 6).Value = 100112023
$ws.Cells.Item(137, 7).Value = "Brócoli"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 2000
$ws.Cells.Item(137, 11).Value = 700
$ws.Cells.Item(137, 12).Value = 800
$ws.Cells.Item(137, 13).Value = 750
$ws.Cells.Item(137, 14).Value = "$/unidad"
$ws.Cells.Item(137, 15).Value = "Región Metropolitana"
$ws.Cells.Item(137, 16).Value = 750
$ws.Cells.Item(137, 17).Value = 1
$ws.Cells.Item(137, 18).Value = "Hortaliza"

# --- Append new row 138 (copy of prior row136 data: "Segunda", Volumen 1000) ---
$ws.Cells.Item(138, 1).Value = 11
$ws.Cells.Item(138, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(138, 3).Value = "Bíobío"
$ws.Cells.Item(138, 4).Value = 44400
$ws.Cells.Item(138, 5).Value = 8
$ws.Cells.Item(138, 6).Value = 100112023
$ws.Cells.Item(138, 7).Value = "Brócoli"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Segunda"
$ws.Cells.Item(138, 10).Value = 1000
$ws.Cells.Item(138, 11).Value = 600
$ws.Cells.Item(138, 12).Value = 600
$ws.Cells.Item(138, 13).Value = 600
$ws.Cells.Item(138, 14).Value = "$/unidad"
$ws.Cells.Item(138, 15).Value = "Región Metropolitana"
$ws.Cells.Item(138, 16).Value = 600
$ws.Cells.Item(138, 17).Value = 1
$ws.Cells.Item(138, 18).Value = "Hortaliza"

# Match the date-cell format used for existing Fecha column (column D) cells.
$ws.Range("D137").NumberFormat = $ws.Range("D136").NumberFormat
$ws.Range("D138").NumberFormat = $ws.Range("D136").NumberFormat
